$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 18: C7 / super capacitor -------------------------------------
# Clone row 17's per-column styling (s="2" ... s="3" for M) by copying
# the fully-populated row 17 and inserting the copy as the new row 18.
$ws.Rows(17).Copy()
$ws.Rows(18).Insert()

$ws.Range("A18").Value = "C7"
$ws.Range("B18").Value = 1
$ws.Range("C18").Clear()
$ws.Range("D18").Value = ""
$ws.Range("E18").Value = ""
$ws.Range("F18").Value = "super capacitor"
$ws.Range("G18").Value = ""
$ws.Range("H18").Value = "top"
$ws.Range("I18").Value = "y"
$ws.Range("J18").Value = "JUWT1105MCD"
$ws.Range("K18").Value = "2.7V"
$ws.Range("L18").Value = ""
$ws.Range("M18").Value = 0.2
$ws.Range("N18").Value = "Nichicon"
$ws.Range("O18").Value = "JUWT1105MCD"
$ws.Range("P18").Value = ""

# --- Row 19: P1 / solar panel -------------------------------------
$ws.Rows(17).Copy()
$ws.Rows(19).Insert()

$ws.Range("A19").Value = "P1"
$ws.Range("B19").Value = 1
$ws.Range("C19").Clear()
$ws.Range("D19").Value = ""
$ws.Range("E19").Value = ""
$ws.Range("F19").Value = "solar panel"
$ws.Range("G19").Value = ""
$ws.Range("H19").Clear()
$ws.Range("I19").Value = "n"
$ws.Range("J19").Value = "1V 80mA 30x25mm"
$ws.Range("K19").Value = "1V"
$ws.Range("L19").Value = ""
$ws.Range("M19").ClearContents()
$ws.Range("N19").Value = "Generic part"
$ws.Range("O19").Value = "Generic part"
$ws.Range("P19").Value = ""

$ws.Range("L19").Select()
